$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank "Nov" (column L) values, and correct a few
# "Sep"/"Oct" (columns J/K) values that shifted as later weekly/monthly
# survey data came in (EIA-914 Feb-2017 snapshot).

# Row 6
$ws.Range("L6").Value = 23

# Row 7
$ws.Range("L7").Value = 513

# Row 8
$ws.Range("L8").Value = 0

# Row 9
$ws.Range("L9").Value = 14

# Row 10
$ws.Range("L10").Value = 505

# Row 11
$ws.Range("L11").Value = 307

# Row 12
$ws.Range("L12").Value = 5

# Row 13
$ws.Range("L13").Value = 23

# Row 14
$ws.Range("L14").Value = 5

# Row 15
$ws.Range("L15").Value = 95

# Row 16
$ws.Range("L16").Value = 6

# Row 17
$ws.Range("L17").Value = 145

# Row 18
$ws.Range("L18").Value = 15

# Row 19
$ws.Range("J19").Value = 57
$ws.Range("K19").Value = 56
$ws.Range("L19").Value = 55

# Row 20
$ws.Range("L20").Value = 0

# Row 21
$ws.Range("L21").Value = 61

# Row 22
$ws.Range("L22").Value = 6

# Row 23
$ws.Range("L23").Value = 1

# Row 24
$ws.Range("J24").Value = 402
$ws.Range("K24").Value = 409
$ws.Range("L24").Value = 417

# Row 25
$ws.Range("L25").Value = 1

# Row 26
$ws.Range("L26").Value = 1029

# Row 27
$ws.Range("L27").Value = 48

# Row 28
$ws.Range("J28").Value = 411
$ws.Range("K28").Value = 422
$ws.Range("L28").Value = 422

# Row 29
$ws.Range("L29").Value = 16

# Row 30
$ws.Range("L30").Value = 4

# Row 31
$ws.Range("L31").Value = 1

# Row 32
$ws.Range("J32").Value = 3152
$ws.Range("K32").Value = 3175
$ws.Range("L32").Value = 3195

# Row 33
$ws.Range("L33").Value = 83

# Row 34
$ws.Range("L34").Value = 0

# Row 35
$ws.Range("L35").Value = 21

# Row 36
$ws.Range("J36").Value = 186
$ws.Range("K36").Value = 192
$ws.Range("L36").Value = 192

# Row 38
$ws.Range("K38").Value = 1591
$ws.Range("L38").Value = 1681

# Row 39
$ws.Range("J39").Value = 16
$ws.Range("K39").Value = 16
$ws.Range("L39").Value = 14

# Row 40
$ws.Range("J40").Value = 8567
$ws.Range("K40").Value = 8799
$ws.Range("L40").Value = 8904

# Move the active cell selection on the sheet to A3 (was T24).
$ws.Range("A3").Select()
